$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------------
# 1) Insert a brand-new row at position 3 for "A 34759-2023" (it moves from
#    its old position - row 6 - up to row 3, with refreshed figures), which
#    shifts the previous rows 3,4,5,6 down to 4,5,6,7.
# -----------------------------------------------------------------------
$ws.Rows("3:3").Insert()

# -----------------------------------------------------------------------
# 2) The old "A 34759-2023" row (now sitting at row 7 after the shift above)
#    is now a duplicate of the new row 3, so remove it - this brings every
#    later row back to its original row number.
# -----------------------------------------------------------------------
$ws.Rows("7:7").Delete()

# -----------------------------------------------------------------------
# 3) Fill the new row 3 with the refreshed "A 34759-2023" data.
# -----------------------------------------------------------------------
$ws.Rows("3:3").RowHeight = 15

$ws.Range("A3").Value = "A 34759-2023"
$ws.Range("B3").Value = 45139
$ws.Range("C3").Value = 45203
$ws.Range("D3").Value = "VÄSTMANLANDS LÄN"
$ws.Range("E3").Value = "SURAHAMMAR"
$ws.Range("F3").Value = "Bergvik skog väst AB"
$ws.Range("G3").Value = 38.3
$ws.Range("H3").Value = 7
$ws.Range("I3").Value = 5
$ws.Range("J3").Value = 7
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 8
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 16
$ws.Range("R3").Value = "Knärot`r`nBlå taggsvamp`r`nGrantaggsvamp`r`nMotaggsvamp`r`nSkogshare`r`nSpillkråka`r`nTalltita`r`nUllticka`r`nDropptaggsvamp`r`nGrönpyrola`r`nMindre märgborre`r`nPlattlummer`r`nVedticka`r`nLopplummer`r`nMattlummer`r`nRevlummer"

$ws.Range("S3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_SURAHAMMAR/artfynd/A 34759-2023.xlsx", "A 34759-2023")'
$ws.Range("T3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_SURAHAMMAR/kartor/A 34759-2023.png", "A 34759-2023")'
$ws.Range("U3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_SURAHAMMAR/knärot/A 34759-2023.png", "A 34759-2023")'
$ws.Range("V3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_SURAHAMMAR/klagomål/A 34759-2023.docx", "A 34759-2023")'
$ws.Range("W3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_SURAHAMMAR/klagomålsmail/A 34759-2023.docx", "A 34759-2023")'
$ws.Range("X3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_SURAHAMMAR/tillsyn/A 34759-2023.docx", "A 34759-2023")'
$ws.Range("Y3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_SURAHAMMAR/tillsynsmail/A 34759-2023.docx", "A 34759-2023")'

# -----------------------------------------------------------------------
# 4) Every data row's "Förändrad" date (column C) moves from 45202 to
#    45203 - one day later. Apply this across the whole data range in one
#    shot (row 3, just populated above, already carries 45203).
# -----------------------------------------------------------------------
$ws.Range("C2:C199").Value = 45203
